$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# REP-13: Move name in template
# Split the combined "position + full name" labels (both the literal Russian
# caption in column A and the corresponding placeholder in column E) across
# two rows instead of one.

# Row 41 previously held "Председатель Правления Бурлаков К.В." - now only
# the position remains here.
$ws.Range("A41").Value = "Председатель Правления"

# Row 42 previously was empty - it now receives the name that used to share
# row 41 with the position.
$ws.Range("A42").Value = "Бурлаков К.В."

# Same split for the template placeholders in column E.
# (single-quoted so PowerShell does not try to expand the leading '$'.)
$ws.Range("E41").Value = '${partyRepresentation.merchantRepresentativePosition}'
$ws.Range("E42").Value = '${partyRepresentation.merchantRepresentativeFullName}'

# The active selection moved from E10 to H30 in the saved workbook.
$ws.Range("H30").Select()
